$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 00:52:04"
$wsZhCn.Range("H2").Value = "2016-03-20 00:52:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 00:52:07"
$wsDeDe.Range("H2").Value = "2016-03-20 00:52:28"
